# Update column C (Förändrad / date changed) from 45243 (2023-11-13) to
# 45244 (2023-11-14) for rows 2 through 10 on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 10; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45243) {
        $cell.Value2 = 45244
    }
}
